$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (8) to make room for "Augment"
$ws.Columns.Item(8).Insert()

# Set the new header cell H1 = "Augment" (inherits the bold/centered header style
# automatically from the column-insert operation, matching the other header cells)
$ws.Cells.Item(1, 8).Value = "Augment"

# Add the five new Artificer skill rows (182-186)
# Row 182: Appraise
$ws.Cells.Item(182, 1).Value = 'Appraise'
$ws.Cells.Item(182, 2).Value = 'Upon taking the Appraise skill, an artificer can select ONE expertise to be well-versed in (non-exhaustive examples in parentheses). Gems (Bangles, Broaches, Rings), Runes (Tattoos, Brands, Etchings), Weaponry (Blades, Hammers, Spears, etc.), Armor (Breastplates, Greaves, Helms, Shields, etc.), Spell Focus (Wizard Staves, Spellbooks, Holy Symbols, etc.), Alchemy (Potions, Reagents, Monster Parts, etc.), Machinery (Clockwork Machines, Engines, Traps, etc.), Clothing (Boots, Capes, Cloaks, Hats, etc.),Art (Paintings, Sculptures, Tapestries, etc.). The Artificer also gains the ability to write Schematics to replicate the effects of spells– embedding in Artificer Objects of your Appraisal expertise.  The Artificer may read scrolls, magical runes, and other materials as though they had the Read/Write Arcana spell, but only to understand them as a reference for creating items (see The Object’s Schematics below), not to cast the spells.  When scrolls are read only as a reference, they are neither activated nor expended. When using them to create a Schematic, the scroll or other reference material is expended due to experimentation.  The Schematic persists indefinitely for future projects, but must be at-hand to create an Artificer Object (they can be stolen through Pickpocket or Loot and are prioritized over coin).'
$ws.Cells.Item(182, 3).Value = 'Artificer'
$ws.Cells.Item(182, 4).Value = 1
$ws.Cells.Item(182, 5).Value = 'In terms of roleplay, your keen eye for objects of this category allows you to intuit its relative worth and even altered properties, subject to Organizer oversight.  Make it fun and specific to your character!  You may take this skill a number of times equal to your Artificer Tier, each time selecting a new Expertise.  An Artificer can only understand reference materials in Tier up to their Artificer Professional Expertise level. An Artificer may read as many reference materials as they like, but may only draft one Schematic per event.  These Schematics are mundane and may not be used to cast spells like scrolls.  Once drafted, a Schematic does not scale—even if an Artificer that authored it improves, the Schematic stays the same level.'
$ws.Cells.Item(182, 6).Value = 'To Appraise - a magnifying glass, examiner’s loupe, or craftsperson-specific object.  For Schematics - an artistic architectural document that shows the fashioning of an object or demonstration of an ability.'
$ws.Cells.Item(182, 9).Value = $False

# Row 183: Tinkerer’s Quirk
$ws.Cells.Item(183, 1).Value = 'Tinkerer’s Quirk'
$ws.Cells.Item(183, 2).Value = 'You may attach a mechanism to an object that alters its properties in one way of your choice with no gem cost - the object (1) sheds light, as though through a light spell, (2) shows a static visual effect (or small picture), (3) makes a continuous non-damaging audio effect or a nonverbal sound, or (4) plays a recorded message when touched, up to 6 seconds long.'
$ws.Cells.Item(183, 3).Value = 'Artificer'
$ws.Cells.Item(183, 4).Value = 1
$ws.Cells.Item(183, 6).Value = 'an actual light emitting material, a picture, a small speaker, or a message strip as appropriate adhered to or within another object.  Make it!'
$ws.Cells.Item(183, 9).Value = $False

# Row 184: Schematic Encryption
$ws.Cells.Item(184, 1).Value = 'Schematic Encryption'
$ws.Cells.Item(184, 2).Value = 'You may disguise the meaning of your Artificer Schematics, so other fellow Artificers cannot read it to learn its secrets.  Upon creation of a Schematic, you may opt to encrypt it.  The Professional expertise required to bypass your encryption for this Schematic is your Artificer Professional Expertise Tier when the Schematic was created.'
$ws.Cells.Item(184, 3).Value = 'Artificer'
$ws.Cells.Item(184, 4).Value = 1
$ws.Cells.Item(184, 6).Value = 'A Schematic (see The Object’s Schematics below) that hides the description, tier, and other information behind a card or other device that is labeled with the Artificer Tier necessary to decrypt.'
$ws.Cells.Item(184, 9).Value = $False

# Row 185: Artificer’s Cooperative Action
$ws.Cells.Item(185, 1).Value = 'Artificer’s Cooperative Action'
$ws.Cells.Item(185, 2).Value = 'The best Artificers learn that they can get further by working together. When faced with a Schematic that is beyond their capabilities (Tier), they can utilize Cooperative Action to encrypt or decrypt Schematics and gain an advantage. The Artificer that takes the action receives +1 to their effective Professional experience Tier for each Artificer with Cooperative Action working together. In order for the Cooperative Action to be successful, the total tier must be greater than the encryption to succeed. It takes double the amount of time to complete a Cooperative Action as it would to complete the action normally.'
$ws.Cells.Item(185, 3).Value = 'Artificer'
$ws.Cells.Item(185, 4).Value = 1
$ws.Cells.Item(185, 5).Value = 'Keep in mind that if you encrypt a Schematic at a Tier higher than you are capable of that you will be unable to utilize it as a reference material without the help of others!'
$ws.Cells.Item(185, 6).Value = 'All Artificers must be within touch distance of the Schematic and actively engaged in the roleplay.'
$ws.Cells.Item(185, 7).Value = 'Schematic Encryption'
$ws.Cells.Item(185, 9).Value = $False

# Row 186: Artorias’ Silver Weapon
$ws.Cells.Item(186, 1).Value = 'Artorias’ Silver Weapon'
$ws.Cells.Item(186, 2).Value = 'Through careful art, science, and ritual, you are able to silver-plate to your weapon’s blade like Sir Artorias the Moonslayer.  Until the end of the day, your blade is silvered, allowing your attacks with this weapon to be harrowing against lycanthropes, wraiths, vampires, and other cursed entities.'
$ws.Cells.Item(186, 3).Value = 'Artificer'
$ws.Cells.Item(186, 4).Value = 1
$ws.Cells.Item(186, 6).Value = 'A bladed weapon.  A number of silver pieces that, laid flat against the blade portion to cover, multiplied by three, is consumed in the act of silvering it (returned to Organizers).  The bigger the blade, the more surface area to cover in coin.  The silver coin, once applied, is consumed.  Once complete, applying a black ribbon with three silver stripes to the weapon indicates that it is silvered.'
$ws.Cells.Item(186, 8).Value = 'Expending a Small gem and a Fusion Point, ground in during the alchemical process, instead renders the weapon silvered for an entire event.  Expending a Large gem and a Fusion Point during its creation renders the silver plating permanent on the weapon.'
$ws.Cells.Item(186, 9).Value = $False
